$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for updated coin rows.
# Force text number format on each touched cell first so numeric-looking
# strings (with significant trailing zeros, % signs, leading zeros, etc.)
# are stored verbatim as text instead of being parsed into floating point
# numbers (which would silently drop precision, e.g. "0.1400" -> 0.14).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "260.56"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.14%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.28"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.37%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.684"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.89%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06114"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.21%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.656"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.39%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8519"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.07%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9217"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.68%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1400"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.71%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04762"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "14.35%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07086"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.46%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03052"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.67%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09065"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.29%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001531"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.23%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006070"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.69%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006046"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.15%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.451"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.55%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.14%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.164"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.60%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.92%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1310"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.02%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.110"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "6.56%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04224"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.41%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001221"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.43%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.003804"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.38%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03854"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.28%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1112"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.73%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004091"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-34.52%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "12.47%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.79%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005167"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.28%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.04%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "35.03%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1593"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-33.88%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.04%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"
